$d = $word.ActiveDocument

function Replace-Exact($oldText, $newText) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -ge 0) {
        $rng = $d.Range($idx, $idx + $oldText.Length)
        $rng.Text = $newText
    }
}

Replace-Exact "Username: joe@gmail.com" "Username: ronald@gmail.com"
Replace-Exact "Password: password123" "Password: password"
Replace-Exact "Firstname: joe" "Firstname: ronald"
Replace-Exact "Surname: swans" "Surname: noble"
Replace-Exact "Telephone: 47583629936" "Telephone: 07758607064"
Replace-Exact "Postcode: y" "Postcode: bt26 "
Replace-Exact "Age: 34" "Age: 43"
Replace-Exact "Group: 7" "Group: 8"
